$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.939.02"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.634.72"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.77"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.69"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "1.860.70"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "1.616.03"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.90"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0₃0755"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "25.928.49"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.25"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.94"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  +2.05%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.49"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "1.138.62"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.804"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.46"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.18"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").Value = "1.770.03"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.34"
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0525"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.60"
$ws.Range("E51").Value = "  +1.36%  "
